# Append the 2025-01-25 00:21 resale-number snapshot as new row 41.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 41

# Columns A-D are text in this sheet (Date/Time/Weekday/Week-as-zero-padded
# string). Force text interpretation so "2025-01-25"/"03" aren't coerced to
# a date serial / number, then drop the resulting number-format override so
# the new row's cells stay styleless like the rest of the data rows.
$textRange = $ws.Range(("A" + $row + ":D" + $row))
$textRange.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-01-25"
$ws.Range("B" + $row).Value = "00:21:41"
$ws.Range("C" + $row).Value = "Saturday"
$ws.Range("D" + $row).Value = "03"

$textRange.ClearFormats()

# Columns E-T are the per-city numeric resale counts (-1 = no data).
$ws.Range("E" + $row).Value = 126181
$ws.Range("F" + $row).Value = 142087
$ws.Range("G" + $row).Value = 168463
$ws.Range("H" + $row).Value = 158589
$ws.Range("I" + $row).Value = -1
$ws.Range("J" + $row).Value = 142794
$ws.Range("K" + $row).Value = -1
$ws.Range("L" + $row).Value = -1
$ws.Range("M" + $row).Value = 191598
$ws.Range("N" + $row).Value = 115686
$ws.Range("O" + $row).Value = 45560
$ws.Range("P" + $row).Value = 28391
$ws.Range("Q" + $row).Value = 65269
$ws.Range("R" + $row).Value = -1
$ws.Range("S" + $row).Value = 47112
$ws.Range("T" + $row).Value = -1
